# "Delete exam and batch info done!"
#
# - Fills in the Maths/Science maximum-marks row (H7/I7).
# - Fills in the first three student rows (8-10) with name, subject marks,
#   a running Total formula, the already-entered max total, and a remark.
# - Clears the left-over serial numbers that had been pre-filled for the
#   still-empty rows 11-20 (batch info no longer needed there), restoring
#   the plain unbordered-top style those rows should carry.
# - Leaves the active selection on D9, matching the saved workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Maximum marks row (row 7) - Maths (H) / Science (I) maximums.
$ws.Cells.Item(7, 8).Value = 70
$ws.Cells.Item(7, 9).Value = 50

# --- Names (column B) first, so new shared strings are interned in the
#     same order the source workbook used. ---
$ws.Cells.Item(8, 2).Value = "Deb"
$ws.Cells.Item(9, 2).Value = "Gupta"
$ws.Cells.Item(10, 2).Value = "Sharma"

# --- Subject marks + totals for rows 8-10. ---
$ws.Cells.Item(8, 3).Value = 19
$ws.Cells.Item(8, 4).Value = 5
$ws.Cells.Item(8, 5).Value = 6
$ws.Cells.Item(8, 6).Value = 10
$ws.Cells.Item(8, 7).Value = 15
$ws.Cells.Item(8, 8).Formula = "=C8+D8+E8+F8+G8"
$ws.Cells.Item(8, 9).Value = 25

$ws.Cells.Item(9, 3).Value = 15
$ws.Cells.Item(9, 4).Value = 5
$ws.Cells.Item(9, 5).Value = 10
$ws.Cells.Item(9, 6).Value = 9
$ws.Cells.Item(9, 7).Value = 14
$ws.Cells.Item(9, 9).Value = 48

$ws.Cells.Item(10, 3).Value = 13
$ws.Cells.Item(10, 4).Value = 7
$ws.Cells.Item(10, 5).Value = 9
$ws.Cells.Item(10, 6).Value = 8
$ws.Cells.Item(10, 7).Value = 18
$ws.Cells.Item(10, 9).Value = 50

# H9/H10 pick up the same "box" formatting H8 already has (loses the thin
# top border that used to separate each still-empty row) before getting
# their running-total formula.
$ws.Cells.Item(8, 8).Copy()
$ws.Cells.Item(9, 8).PasteSpecial(-4122)
$ws.Cells.Item(10, 8).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(9, 8).Formula = "=C9+D9+E9+F9+G9"
$ws.Cells.Item(10, 8).Formula = "=C10+D10+E10+F10+G10"

# --- Remarks (column J) last. ---
$ws.Cells.Item(8, 10).Value = "Legen-wait for it-dary"
$ws.Cells.Item(9, 10).Value = "Pathetic"
$ws.Cells.Item(10, 10).Value = "Awesome"

# Clear the serial numbers pre-filled in rows 11-20 (no batch assigned to
# those rows yet) and restore the unbordered-top look used by the rest of
# the still-empty rows, by pulling formatting from an already-plain row.
for ($r = 11; $r -le 20; $r++) {
    $ws.Cells.Item(21, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
    $ws.Cells.Item($r, 1).ClearContents()
}
$excel.CutCopyMode = $false

# Move the active selection, as recorded in the saved workbook.
$ws.Range("D9").Select()
